# Add more blank slides to default presentation
#
# The source deck ends with a slide whose only content is a text box
# reading "This slide left blank for whiteboard". The author's edit
# appends 14 more copies of that same slide to the end of the deck
# (going from 6 slides to 20), so there are more blank "whiteboard"
# slides available. No other slide content changes.

$p = $ppt.ActivePresentation

$blankSlideCount = 14

for ($i = 0; $i -lt $blankSlideCount; $i++) {
    $lastSlide = $p.Slides.Item($p.Slides.Count)
    $lastSlide.Duplicate() | Out-Null
}

Write-Output "Slides now: $($p.Slides.Count)"
